# Fruta / hortaliza, semanal
# Re-shuffle the D, J, K, L, M, O, P values across rows 2-13 according to the
# permutation observed between the "before" and "after" versions of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (values to copy FROM source row INTO target row)
$mapping = @{
    2  = 13
    3  = 4
    4  = 9
    5  = 10
    6  = 7
    7  = 6
    8  = 3
    9  = 12
    10 = 5
    11 = 2
    12 = 8
    13 = 11
}

# Columns that are permuted together as a record.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot all current (pre-edit) values for the affected columns/rows so that
# later writes don't clobber values still needed as a source for other rows.
# Value2 is used (rather than Value) because it reliably returns the raw
# cell contents (numbers/strings) in this environment.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Apply the permutation.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceValues = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $sourceValues[$col]
    }
}
